$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.905.66"
$ws.Range("E2").Value = "  -1.02%  "
$ws.Range("D3").Value = "2.900.31"
$ws.Range("E3").Value = "  -1.50%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "568.82"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -3.53%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.86"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.04%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  -1.19%  "
$ws.Range("D9").Value = "2.897.45"
$ws.Range("E9").Value = "  -1.55%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.91"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -2.74%  "
$ws.Range("E11").Value = "  -3.50%  "
$ws.Range("E12").Value = "  -1.15%  "
$ws.Range("E13").Value = "  -1.34%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.17"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.37%  "
$ws.Range("E15").Value = "  -0.19%  "
$ws.Range("D16").Value = "3.380.98"
$ws.Range("E16").Value = "  -1.50%  "
$ws.Range("D17").Value = "61.851.13"
$ws.Range("E17").Value = "  -1.01%  "
$ws.Range("D18").Value = "2.897.04"
$ws.Range("E18").Value = "  -1.60%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.50"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -2.34%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "431.38"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.76%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.92"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -4.17%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.652"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.74%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.86"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.76%  "
$ws.Range("E24").Value = "  -1.70%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.97"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.62%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.10"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -9.41%  "
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.02"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -3.77%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0000111"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +9.47%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.01"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -2.62%  "
$ws.Range("E31").Value = "  -2.94%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.02"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -6.26%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.07%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.106"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -2.24%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "25.61"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -2.29%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.955"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -3.47%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.39"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -3.43%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "48.80"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.79%  "
$ws.Range("E39").Value = "  -5.55%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.91"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -4.99%  "
$ws.Range("E41").Value = "  -0.34%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.14"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -2.73%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "40.33"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +4.48%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.268"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -2.88%  "
$ws.Range("D45").Value = "2.700.75"
$ws.Range("E45").Value = "  +0.36%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0337"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.49%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "131.63"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.67%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "346.21"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -2.39%  "
$ws.Range("E50").Value = "  -1.14%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "21.58"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -4.31%  "
